$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ParentID of row 4 ("Haben Sie weitere Fahrzeuge?") from 2 to 1
$ws.Range("B4").Value = 1

# Move the active selection to B5 (also clears the scrolled topLeftCell)
$ws.Range("B5").Select()
